$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing data right
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy formatting from column F (old column D, now shifted) into new D:E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(7, 4).Value = 43434
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = "NA"
$ws.Cells.Item(8, 5).Value = 1459000
$ws.Cells.Item(8, 6).Value = 2404900
$ws.Cells.Item(8, 7).Value = 1161100
$ws.Cells.Item(8, 8).Value = 3194100
$ws.Cells.Item(8, 9).Value = 1102600
$ws.Cells.Item(8, 10).Value = 2633100
$ws.Cells.Item(8, 11).Value = 1518900
$ws.Cells.Item(8, 12).Value = 2921600
$ws.Cells.Item(8, 13).Value = 2883700
$ws.Cells.Item(9, 4).Value = 441200
$ws.Cells.Item(9, 5).Value = 437600
$ws.Cells.Item(9, 6).Value = 858200
$ws.Cells.Item(9, 7).Value = 389800
$ws.Cells.Item(9, 8).Value = 2101000
$ws.Cells.Item(9, 9).Value = 359800
$ws.Cells.Item(9, 10).Value = 699300
$ws.Cells.Item(9, 11).Value = 281600
$ws.Cells.Item(9, 12).Value = 1934400
$ws.Cells.Item(9, 13).Value = 1896400
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(10, 5).Value = 1021400
$ws.Cells.Item(10, 6).Value = 1546700
$ws.Cells.Item(10, 7).Value = 771300
$ws.Cells.Item(10, 8).Value = 1093100
$ws.Cells.Item(10, 9).Value = 742800
$ws.Cells.Item(10, 10).Value = 1933800
$ws.Cells.Item(10, 11).Value = 1237300
$ws.Cells.Item(10, 12).Value = 987200
$ws.Cells.Item(10, 13).Value = 987300
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(15, 4).Value = 28000
$ws.Cells.Item(15, 5).Value = 32300
$ws.Cells.Item(15, 6).Value = 60100
$ws.Cells.Item(15, 7).Value = 28200
$ws.Cells.Item(15, 8).Value = 53300
$ws.Cells.Item(15, 9).Value = 28800
$ws.Cells.Item(15, 10).Value = 53400
$ws.Cells.Item(15, 11).Value = 27100
$ws.Cells.Item(15, 12).Value = 58500
$ws.Cells.Item(15, 13).Value = 52700
$ws.Cells.Item(17, 4).Value = 915500
$ws.Cells.Item(17, 5).Value = 1205200
$ws.Cells.Item(17, 6).Value = 2395100
$ws.Cells.Item(17, 7).Value = 1155400
$ws.Cells.Item(17, 8).Value = 2914400
$ws.Cells.Item(17, 9).Value = 1107100
$ws.Cells.Item(17, 10).Value = 2131400
$ws.Cells.Item(17, 11).Value = 1023800
$ws.Cells.Item(17, 12).Value = 2719500
$ws.Cells.Item(17, 13).Value = 2647200
$ws.Cells.Item(18, 4).Value = "NA"
$ws.Cells.Item(18, 5).Value = 253800
$ws.Cells.Item(18, 6).Value = 9800
$ws.Cells.Item(18, 7).Value = 5700
$ws.Cells.Item(18, 8).Value = 279700
$ws.Cells.Item(18, 9).Value = -4500
$ws.Cells.Item(18, 10).Value = 501700
$ws.Cells.Item(18, 11).Value = 495100
$ws.Cells.Item(18, 12).Value = 202100
$ws.Cells.Item(18, 13).Value = 236500
$ws.Cells.Item(20, 4).Value = "NA"
$ws.Cells.Item(20, 5).Value = 18800
$ws.Cells.Item(20, 6).Value = 65500
$ws.Cells.Item(20, 7).Value = 32100
$ws.Cells.Item(20, 8).Value = 11300
$ws.Cells.Item(20, 9).Value = 30100
$ws.Cells.Item(20, 10).Value = -114500
$ws.Cells.Item(20, 11).Value = -128600
$ws.Cells.Item(20, 12).Value = 67500
$ws.Cells.Item(20, 13).Value = 36500
$ws.Cells.Item(21, 4).Value = "NA"
$ws.Cells.Item(21, 5).Value = 300800
$ws.Cells.Item(21, 6).Value = 127800
$ws.Cells.Item(21, 7).Value = 76500
$ws.Cells.Item(21, 8).Value = 367700
$ws.Cells.Item(21, 9).Value = 50300
$ws.Cells.Item(21, 10).Value = 432000
$ws.Cells.Item(21, 11).Value = 404600
$ws.Cells.Item(21, 12).Value = 315900
$ws.Cells.Item(21, 13).Value = 315200
$ws.Cells.Item(22, 4).Value = "NA"
$ws.Cells.Item(22, 5).Value = "NA"
$ws.Cells.Item(22, 6).Value = "NA"
$ws.Cells.Item(22, 7).Value = "NA"
$ws.Cells.Item(22, 8).Value = "NA"
$ws.Cells.Item(22, 9).Value = "NA"
$ws.Cells.Item(22, 10).Value = "NA"
$ws.Cells.Item(22, 11).Value = 25600
$ws.Cells.Item(22, 12).Value = 40600
$ws.Cells.Item(22, 13).Value = 23100
$ws.Cells.Item(23, 4).Value = "NA"
$ws.Cells.Item(23, 5).Value = 272700
$ws.Cells.Item(23, 6).Value = 75300
$ws.Cells.Item(23, 7).Value = 37800
$ws.Cells.Item(23, 8).Value = 291000
$ws.Cells.Item(23, 9).Value = 25500
$ws.Cells.Item(23, 10).Value = 387300
$ws.Cells.Item(23, 11).Value = 341000
$ws.Cells.Item(23, 12).Value = 229000
$ws.Cells.Item(23, 13).Value = 249900
$ws.Cells.Item(24, 4).Value = -32600
$ws.Cells.Item(24, 5).Value = 90400
$ws.Cells.Item(24, 6).Value = -38800
$ws.Cells.Item(24, 7).Value = -48400
$ws.Cells.Item(24, 8).Value = 92400
$ws.Cells.Item(24, 9).Value = 9800
$ws.Cells.Item(24, 10).Value = 117400
$ws.Cells.Item(24, 11).Value = 91200
$ws.Cells.Item(24, 12).Value = 62900
$ws.Cells.Item(24, 13).Value = 73700
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = "NA"
$ws.Cells.Item(26, 5).Value = 182300
$ws.Cells.Item(26, 6).Value = 114100
$ws.Cells.Item(26, 7).Value = 86200
$ws.Cells.Item(26, 8).Value = 198600
$ws.Cells.Item(26, 9).Value = 15800
$ws.Cells.Item(26, 10).Value = 269800
$ws.Cells.Item(26, 11).Value = 249800
$ws.Cells.Item(26, 12).Value = 166100
$ws.Cells.Item(26, 13).Value = 176200
$ws.Cells.Item(27, 4).Value = "NA"
$ws.Cells.Item(27, 5).Value = 191600
$ws.Cells.Item(27, 6).Value = 72100
$ws.Cells.Item(27, 7).Value = 71600
$ws.Cells.Item(27, 8).Value = 180000
$ws.Cells.Item(27, 9).Value = -22000
$ws.Cells.Item(27, 10).Value = 240100
$ws.Cells.Item(27, 11).Value = 236100
$ws.Cells.Item(27, 12).Value = 136600
$ws.Cells.Item(27, 13).Value = 152300
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 774000
$ws.Cells.Item(29, 7).Value = 53000
$ws.Cells.Item(29, 8).Value = -450500
$ws.Cells.Item(29, 9).Value = 121000
$ws.Cells.Item(29, 10).Value = 152200
$ws.Cells.Item(29, 11).Value = 44200
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = "NA"
$ws.Cells.Item(32, 5).Value = -18800
$ws.Cells.Item(32, 6).Value = -65500
$ws.Cells.Item(32, 7).Value = -32100
$ws.Cells.Item(32, 8).Value = -11300
$ws.Cells.Item(32, 9).Value = -30100
$ws.Cells.Item(32, 10).Value = 114500
$ws.Cells.Item(32, 11).Value = 128600
$ws.Cells.Item(32, 12).Value = -67500
$ws.Cells.Item(32, 13).Value = -36500
$ws.Cells.Item(33, 4).Value = "NA"
$ws.Cells.Item(33, 5).Value = 191600
$ws.Cells.Item(33, 6).Value = 846100
$ws.Cells.Item(33, 7).Value = 124600
$ws.Cells.Item(33, 8).Value = -270500
$ws.Cells.Item(33, 9).Value = 99000
$ws.Cells.Item(33, 10).Value = 392300
$ws.Cells.Item(33, 11).Value = 280300
$ws.Cells.Item(33, 12).Value = 136600
$ws.Cells.Item(33, 13).Value = 152300
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = "NA"
$ws.Cells.Item(35, 5).Value = 191600
$ws.Cells.Item(35, 6).Value = 846100
$ws.Cells.Item(35, 7).Value = 124600
$ws.Cells.Item(35, 8).Value = -270500
$ws.Cells.Item(35, 9).Value = 99000
$ws.Cells.Item(35, 10).Value = 392300
$ws.Cells.Item(35, 11).Value = 280300
$ws.Cells.Item(35, 12).Value = 136600
$ws.Cells.Item(35, 13).Value = 152300
$ws.Cells.Item(38, 4).Value = 43434
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(41, 4).Value = 5258800
$ws.Cells.Item(41, 5).Value = 4895800
$ws.Cells.Item(41, 6).Value = 4741100
$ws.Cells.Item(41, 7).Value = 5144600
$ws.Cells.Item(41, 8).Value = 5275500
$ws.Cells.Item(41, 9).Value = 5016900
$ws.Cells.Item(41, 10).Value = 4661900
$ws.Cells.Item(41, 11).Value = 4966700
$ws.Cells.Item(41, 12).Value = 4664900
$ws.Cells.Item(41, 13).Value = 4305800
$ws.Cells.Item(42, 4).Value = 9324000
$ws.Cells.Item(42, 5).Value = 11029000
$ws.Cells.Item(42, 6).Value = 11421300
$ws.Cells.Item(42, 7).Value = 10284100
$ws.Cells.Item(42, 8).Value = 11411400
$ws.Cells.Item(42, 9).Value = 11130000
$ws.Cells.Item(42, 10).Value = 12245900
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(43, 4).Value = 6287400
$ws.Cells.Item(43, 5).Value = 5864800
$ws.Cells.Item(43, 6).Value = 5876400
$ws.Cells.Item(43, 7).Value = 6710600
$ws.Cells.Item(43, 8).Value = 5419000
$ws.Cells.Item(43, 9).Value = 5507800
$ws.Cells.Item(43, 10).Value = 5874600
$ws.Cells.Item(43, 11).Value = 5553700
$ws.Cells.Item(43, 12).Value = 4425200
$ws.Cells.Item(43, 13).Value = 4169100
$ws.Cells.Item(44, 4).Value = "NA"
$ws.Cells.Item(44, 5).Value = "NA"
$ws.Cells.Item(44, 6).Value = "NA"
$ws.Cells.Item(44, 7).Value = "NA"
$ws.Cells.Item(44, 8).Value = "NA"
$ws.Cells.Item(44, 9).Value = 306300
$ws.Cells.Item(44, 10).Value = 316700
$ws.Cells.Item(44, 11).Value = 310300
$ws.Cells.Item(44, 12).Value = 309900
$ws.Cells.Item(44, 13).Value = 322800
$ws.Cells.Item(45, 4).Value = 708000
$ws.Cells.Item(45, 5).Value = 913500
$ws.Cells.Item(45, 6).Value = 709600
$ws.Cells.Item(45, 7).Value = 802700
$ws.Cells.Item(45, 8).Value = 578000
$ws.Cells.Item(45, 9).Value = 903700
$ws.Cells.Item(45, 10).Value = 919000
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(46, 4).Value = "NA"
$ws.Cells.Item(46, 5).Value = 22703000
$ws.Cells.Item(46, 6).Value = 22748400
$ws.Cells.Item(46, 7).Value = 22942100
$ws.Cells.Item(46, 8).Value = 22683900
$ws.Cells.Item(46, 9).Value = 22864800
$ws.Cells.Item(46, 10).Value = 24018100
$ws.Cells.Item(46, 11).Value = 10830700
$ws.Cells.Item(46, 12).Value = 9400000
$ws.Cells.Item(46, 13).Value = 8797800
$ws.Cells.Item(47, 4).Value = 40163600
$ws.Cells.Item(47, 5).Value = 22028600
$ws.Cells.Item(47, 6).Value = 22084400
$ws.Cells.Item(47, 7).Value = 20270500
$ws.Cells.Item(47, 8).Value = 18866100
$ws.Cells.Item(47, 9).Value = 18917800
$ws.Cells.Item(47, 10).Value = 18327800
$ws.Cells.Item(47, 11).Value = 29044900
$ws.Cells.Item(47, 12).Value = 29532800
$ws.Cells.Item(47, 13).Value = 31108900
$ws.Cells.Item(48, 4).Value = "NA"
$ws.Cells.Item(48, 5).Value = 346900
$ws.Cells.Item(48, 6).Value = 353600
$ws.Cells.Item(48, 7).Value = 754700
$ws.Cells.Item(48, 8).Value = 750400
$ws.Cells.Item(48, 9).Value = 727700
$ws.Cells.Item(48, 10).Value = 728400
$ws.Cells.Item(48, 11).Value = 718400
$ws.Cells.Item(48, 12).Value = 709200
$ws.Cells.Item(48, 13).Value = 715200
$ws.Cells.Item(49, 4).Value = 1890100
$ws.Cells.Item(49, 5).Value = 1894900
$ws.Cells.Item(49, 6).Value = 1901400
$ws.Cells.Item(49, 7).Value = 2451000
$ws.Cells.Item(49, 8).Value = 2463200
$ws.Cells.Item(49, 9).Value = 2473300
$ws.Cells.Item(49, 10).Value = 2488500
$ws.Cells.Item(49, 11).Value = 2499100
$ws.Cells.Item(49, 12).Value = 2513700
$ws.Cells.Item(49, 13).Value = 2594300
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 2372400
$ws.Cells.Item(52, 5).Value = 1975300
$ws.Cells.Item(52, 6).Value = 2314100
$ws.Cells.Item(52, 7).Value = 2583300
$ws.Cells.Item(52, 8).Value = 2405600
$ws.Cells.Item(52, 9).Value = 2683300
$ws.Cells.Item(52, 10).Value = 2813200
$ws.Cells.Item(52, 11).Value = 2886800
$ws.Cells.Item(52, 12).Value = 2915600
$ws.Cells.Item(52, 13).Value = 2889800
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 47131100
$ws.Cells.Item(54, 5).Value = 48948700
$ws.Cells.Item(54, 6).Value = 49401800
$ws.Cells.Item(54, 7).Value = 49001500
$ws.Cells.Item(54, 8).Value = 47169100
$ws.Cells.Item(54, 9).Value = 47666900
$ws.Cells.Item(54, 10).Value = 48376100
$ws.Cells.Item(54, 11).Value = 45980000
$ws.Cells.Item(54, 12).Value = 45071300
$ws.Cells.Item(54, 13).Value = 46106000
$ws.Cells.Item(57, 4).Value = 7407000
$ws.Cells.Item(57, 5).Value = 6680200
$ws.Cells.Item(57, 6).Value = 7496900
$ws.Cells.Item(57, 7).Value = 7336100
$ws.Cells.Item(57, 8).Value = 7167700
$ws.Cells.Item(57, 9).Value = 7105000
$ws.Cells.Item(57, 10).Value = 6785100
$ws.Cells.Item(57, 11).Value = 6859800
$ws.Cells.Item(57, 12).Value = 7373700
$ws.Cells.Item(57, 13).Value = 7765700
$ws.Cells.Item(58, 4).Value = 9030600
$ws.Cells.Item(58, 5).Value = 10246500
$ws.Cells.Item(58, 6).Value = 9279700
$ws.Cells.Item(58, 7).Value = 8718400
$ws.Cells.Item(58, 8).Value = 9096700
$ws.Cells.Item(58, 9).Value = 8890500
$ws.Cells.Item(58, 10).Value = 9060600
$ws.Cells.Item(58, 11).Value = 422900
$ws.Cells.Item(58, 12).Value = 525800
$ws.Cells.Item(58, 13).Value = 432200
$ws.Cells.Item(59, 4).Value = 11317600
$ws.Cells.Item(59, 5).Value = 12010700
$ws.Cells.Item(59, 6).Value = 12738900
$ws.Cells.Item(59, 7).Value = 12364500
$ws.Cells.Item(59, 8).Value = 11298900
$ws.Cells.Item(59, 9).Value = 11531000
$ws.Cells.Item(59, 10).Value = 12569400
$ws.Cells.Item(59, 11).Value = 8764000
$ws.Cells.Item(59, 12).Value = 8388600
$ws.Cells.Item(59, 13).Value = 7999200
$ws.Cells.Item(60, 4).Value = "NA"
$ws.Cells.Item(60, 5).Value = 28937400
$ws.Cells.Item(60, 6).Value = 29515600
$ws.Cells.Item(60, 7).Value = 28418900
$ws.Cells.Item(60, 8).Value = 27563300
$ws.Cells.Item(60, 9).Value = 27526500
$ws.Cells.Item(60, 10).Value = 28415100
$ws.Cells.Item(60, 11).Value = 16046700
$ws.Cells.Item(60, 12).Value = 16288200
$ws.Cells.Item(60, 13).Value = 16197100
$ws.Cells.Item(61, 4).Value = 9276800
$ws.Cells.Item(61, 5).Value = 9343100
$ws.Cells.Item(61, 6).Value = 9304500
$ws.Cells.Item(61, 7).Value = 9879600
$ws.Cells.Item(61, 8).Value = 9040300
$ws.Cells.Item(61, 9).Value = 9110300
$ws.Cells.Item(61, 10).Value = 8972900
$ws.Cells.Item(61, 11).Value = 18998900
$ws.Cells.Item(61, 12).Value = 18142700
$ws.Cells.Item(61, 13).Value = 19325300
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 36925500
$ws.Cells.Item(66, 5).Value = 38327500
$ws.Cells.Item(66, 6).Value = 38863600
$ws.Cells.Item(66, 7).Value = 38742400
$ws.Cells.Item(66, 8).Value = 37063200
$ws.Cells.Item(66, 9).Value = 37215000
$ws.Cells.Item(66, 10).Value = 37884800
$ws.Cells.Item(66, 11).Value = 35581800
$ws.Cells.Item(66, 12).Value = 34943200
$ws.Cells.Item(66, 13).Value = 36046300
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = 5610200
$ws.Cells.Item(72, 5).Value = 5672400
$ws.Cells.Item(72, 6).Value = 5523300
$ws.Cells.Item(72, 7).Value = 4833300
$ws.Cells.Item(72, 8).Value = 4701000
$ws.Cells.Item(72, 9).Value = 5000400
$ws.Cells.Item(72, 10).Value = 4938300
$ws.Cells.Item(72, 11).Value = 4903400
$ws.Cells.Item(72, 12).Value = 4645400
$ws.Cells.Item(72, 13).Value = 4531300
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 10060900
$ws.Cells.Item(76, 5).Value = 10621200
$ws.Cells.Item(76, 6).Value = 10538200
$ws.Cells.Item(76, 7).Value = 10259100
$ws.Cells.Item(76, 8).Value = 10106000
$ws.Cells.Item(76, 9).Value = 10451900
$ws.Cells.Item(76, 10).Value = 10491300
$ws.Cells.Item(76, 11).Value = 10398200
$ws.Cells.Item(76, 12).Value = 10128100
$ws.Cells.Item(76, 13).Value = 10059700
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43434
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = "NA"
$ws.Cells.Item(81, 5).Value = 191600
$ws.Cells.Item(81, 6).Value = 846100
$ws.Cells.Item(81, 7).Value = 124600
$ws.Cells.Item(81, 8).Value = -270500
$ws.Cells.Item(81, 9).Value = 99000
$ws.Cells.Item(81, 10).Value = 392300
$ws.Cells.Item(81, 11).Value = 280300
$ws.Cells.Item(81, 12).Value = 136600
$ws.Cells.Item(81, 13).Value = 152300
$ws.Cells.Item(83, 4).Value = "NA"
$ws.Cells.Item(83, 5).Value = 28100
$ws.Cells.Item(83, 6).Value = 52500
$ws.Cells.Item(83, 7).Value = 38700
$ws.Cells.Item(83, 8).Value = 37100
$ws.Cells.Item(83, 9).Value = 39200
$ws.Cells.Item(83, 10).Value = 69800
$ws.Cells.Item(83, 11).Value = 38000
$ws.Cells.Item(83, 12).Value = 46300
$ws.Cells.Item(83, 13).Value = 42100
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 491200
$ws.Cells.Item(89, 5).Value = 724000
$ws.Cells.Item(89, 6).Value = -524100
$ws.Cells.Item(89, 7).Value = -405300
$ws.Cells.Item(89, 8).Value = 11900
$ws.Cells.Item(89, 9).Value = 587400
$ws.Cells.Item(89, 10).Value = 470400
$ws.Cells.Item(89, 11).Value = -196100
$ws.Cells.Item(89, 12).Value = 272700
$ws.Cells.Item(89, 13).Value = 357100
$ws.Cells.Item(91, 4).Value = -43300
$ws.Cells.Item(91, 5).Value = -41500
$ws.Cells.Item(91, 6).Value = -240900
$ws.Cells.Item(91, 7).Value = -59600
$ws.Cells.Item(91, 8).Value = -103300
$ws.Cells.Item(91, 9).Value = -26600
$ws.Cells.Item(91, 10).Value = -63600
$ws.Cells.Item(91, 11).Value = -94000
$ws.Cells.Item(91, 12).Value = -63300
$ws.Cells.Item(91, 13).Value = -95500
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = 261000
$ws.Cells.Item(94, 5).Value = -126700
$ws.Cells.Item(94, 6).Value = 8100
$ws.Cells.Item(94, 7).Value = -256800
$ws.Cells.Item(94, 8).Value = -25200
$ws.Cells.Item(94, 9).Value = -222900
$ws.Cells.Item(94, 10).Value = 200900
$ws.Cells.Item(94, 11).Value = 125400
$ws.Cells.Item(94, 12).Value = -81100
$ws.Cells.Item(94, 13).Value = -79000
$ws.Cells.Item(96, 4).Value = -40000
$ws.Cells.Item(96, 5).Value = -41700
$ws.Cells.Item(96, 6).Value = -70100
$ws.Cells.Item(96, 7).Value = -36000
$ws.Cells.Item(96, 8).Value = -36000
$ws.Cells.Item(96, 9).Value = -36000
$ws.Cells.Item(96, 10).Value = -45400
$ws.Cells.Item(96, 11).Value = -22700
$ws.Cells.Item(96, 12).Value = -22700
$ws.Cells.Item(96, 13).Value = -22700
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = -538400
$ws.Cells.Item(100, 5).Value = -236300
$ws.Cells.Item(100, 6).Value = 198900
$ws.Cells.Item(100, 7).Value = 874500
$ws.Cells.Item(100, 8).Value = 216800
$ws.Cells.Item(100, 9).Value = -25300
$ws.Cells.Item(100, 10).Value = 243300
$ws.Cells.Item(100, 11).Value = 525700
$ws.Cells.Item(100, 12).Value = 345400
$ws.Cells.Item(100, 13).Value = -18300
$ws.Cells.Item(101, 4).Value = -3100
$ws.Cells.Item(101, 5).Value = -6900
$ws.Cells.Item(101, 6).Value = -9600
$ws.Cells.Item(101, 7).Value = 2600
$ws.Cells.Item(101, 8).Value = 5100
$ws.Cells.Item(101, 9).Value = 700
$ws.Cells.Item(101, 10).Value = 4900
$ws.Cells.Item(101, 11).Value = -500
$ws.Cells.Item(101, 12).Value = -8400
$ws.Cells.Item(101, 13).Value = -11700
$ws.Cells.Item(102, 4).Value = 210700
$ws.Cells.Item(102, 5).Value = 354200
$ws.Cells.Item(102, 6).Value = -326700
$ws.Cells.Item(102, 7).Value = 215100
$ws.Cells.Item(102, 8).Value = 215400
$ws.Cells.Item(102, 9).Value = 338200
$ws.Cells.Item(102, 10).Value = 914300
$ws.Cells.Item(102, 11).Value = 453900
$ws.Cells.Item(102, 12).Value = 528600
$ws.Cells.Item(102, 13).Value = 248000